$d = $word.ActiveDocument

# --- Sanity: locate the anchor paragraph (first empty paragraph right
# after the "...Aproveitamento de caracteristicas ambientais positivas."
# paragraph, and right before the paragraph that carries the _GoBack
# bookmark). We find it by scanning from the end of the document instead
# of hard-coding an index, so the script is resilient if Word's count
# differs slightly.

$total = $d.Paragraphs.Count

# Find the paragraph that holds the (hidden) _GoBack bookmark by walking
# backwards and checking each paragraph range's XML for bookmarkStart.
$bookmarkIndex = -1
for ($i = $total; $i -ge 1; $i--) {
    $xml = $d.Paragraphs($i).Range.WordOpenXML
    if ($xml -like "*_GoBack*") {
        $bookmarkIndex = $i
        break
    }
}
if ($bookmarkIndex -eq -1) {
    throw "Could not locate the _GoBack bookmark paragraph"
}

$anchorIndex = $bookmarkIndex - 1

$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Unidade 3 – Design sustentável</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> (?)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Construção ecológica:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>- uso de resíduos locais;</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>- pequeno aporte de recurso tecnológicos para beneficiamento e transformação de matérias-primas (artesanal). (?)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Ecologiaé a “Relação dos seres vivos com o habitat ou meio ambiente natural” (HAECKEL).O ambiente criado pelo ser humano com o uso de maquinários e tecnologia, </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:lastRenderedPageBreak/><w:t>não se pode falar em relação direta com a natureza, mas com um habitat modificado. Aplicando o mesmo conceito, a definição mais correta para a construção com perfil moderno, urbano e industrial é Construção Sustentável.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Tecnologias sustentáveis-Sistemas ou equipamentos de uso individual, unifamiliar ou ambiente comercial;-Fontes de energia renovável para geração de energia (solar, eólica, biomassa, biodigestores);-Uso, reuso e economia de água;</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Produto ecológico é todo artigo de origem artesanal ou industrializada, de uso pessoal, alimentar, residencial, comercial, agrícola e industrial, que seja não-poluente, não-tóxico, benéfico à saúde e ao meio ambiente, contribuindo para o desenvolvimento de um modelo econômico e social sustentável.Produto sustentável mantém o mesmo conceito + Desempenho sustentável, permitindo o uso racional dos recursos naturais (sem esgotamento).</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Zero, baixa ou controle de geração e emissão de poluentes; gases, efluentes, resíduos sólidos);</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Caráter educativo (comportamento social);</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Matérias primas naturais renováveis ou não renováveis;</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Rótulos Ecológicos (ISO 14024 Tipo I; ISO14025 Tipo II);</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Origem e Controle / DOF (Documento de Origem Florestal); FSC (Forest Stewardship Council); CERFLOR.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Critérios para avaliação</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Matérias-primas empregadas; </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Processo produtivo; </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Consumo de energia elétrica</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Consumo de recursos naturais; </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Poluentes gerados; </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Resíduos gerados; </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Ciclo de vida.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Os modelos oferecem suporte ao projeto ao longo de suas fases.suporte ao projeto ao longo de suas fases.Kam e Fischer (2004) Afirmam que BIM / VDC se baseia na modelagem orientada a processos e a simulações, com o objetivo de avaliar alternativas.Fischer e Kunz (2009) BIM/VDC Simulação digital da integração multidisciplinar dos dados e parâmetros do projeto, contemplando informações para os agentes envolvidos.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Componente Paramétrico = Simulação virtual de um objeto do modelo:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Geometria e dimensões reais </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Representação gráfica dos materiais</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Propriedades físicas (massa, resistência, etc)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Associações </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>a outros elementos construtivos</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:lastRenderedPageBreak/><w:t>Dese</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>mpenho térmico e acústico</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Ciclo de vida</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Garantias </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Código rastreabilidade</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>As “classes” agrupam objetos físicos ou virtuais que apresentem características semelhantes a partir de</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> uma abordagem lógica definida.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>Sistemas: UNIFORMAT, OMNICLASS, ABNT NBR 15965</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-BR" w:bidi="ar-SA"/></w:rPr></w:pPr></w:p>
'@

$anchorRange = $d.Paragraphs($anchorIndex).Range
$anchorRange.InsertXML($fragment)

# Re-resolve the bookmark paragraph (it has shifted down by however many
# paragraphs we just inserted) and drop the (now superfluous) empty
# paragraphs that used to trail it.
$total = $d.Paragraphs.Count
$bookmarkIndex = -1
for ($i = $total; $i -ge 1; $i--) {
    $xml = $d.Paragraphs($i).Range.WordOpenXML
    if ($xml -like "*_GoBack*") {
        $bookmarkIndex = $i
        break
    }
}
if ($bookmarkIndex -eq -1) {
    throw "Could not re-locate the _GoBack bookmark paragraph"
}

# Delete every paragraph that follows the bookmark paragraph (there were
# three empty trailing paragraphs in the source document; the target
# keeps none of them - the bookmark paragraph becomes the very last
# paragraph in the body).
for ($i = $d.Paragraphs.Count; $i -gt $bookmarkIndex; $i--) {
    $d.Paragraphs($i).Range.Delete()
}

Write-Output "paragraphs=$($d.Paragraphs.Count) bookmarkIndex=$bookmarkIndex anchorIndex=$anchorIndex"
